$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Commit: value in C10 (Sample Project / Main.xlsx, sheet "Rules") was
# restored from 18 to 1.
$ws.Range("C10").Value = 1
